$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.023.82"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.62"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.67"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.52"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.76%  "
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("E10").Value = "  -1.78%  "
$ws.Range("E11").Value = "  -3.14%  "
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.60"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.922.33"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.406.06"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "321.20"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("E21").Value = "  +1.75%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.14"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.61"
$ws.Range("E25").Value = "  -5.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "568.33"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0931"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.19"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.08"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.43"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("E39").Value = "  -1.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.54"
$ws.Range("E40").Value = "  -1.12%  "
$ws.Range("E41").Value = "  -3.51%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "148.94"
$ws.Range("E43").Value = "  -1.34%  "
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("E45").Value = "  -0.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0530"
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.594"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.85"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0919"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("E50").Value = "  -0.61%  "
$ws.Range("E51").Value = "  +0.70%  "
